# Update Betfair Back/Lay odds for 2026-02-05 (Sheet1) to match the
# re-scraped/updated values from the source feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - NK Aluminij vs Olimpija
$ws.Range("F2").Value = 4
$ws.Range("H2").Value = 1.66
$ws.Range("K2").Value = 500
$ws.Range("P2").Value = 2.1

# Row 3
$ws.Range("F3").Value = 2.52
$ws.Range("G3").Value = 2.92
$ws.Range("H3").Value = 3.35
$ws.Range("I3").Value = 4.1
$ws.Range("J3").Value = 2.8
$ws.Range("K3").Value = 3.1
$ws.Range("P3").Value = 1.44
$ws.Range("Q3").Value = 2.92

# Row 4
$ws.Range("F4").Value = 1.58
$ws.Range("G4").Value = 1.77
$ws.Range("H4").Value = 5.6
$ws.Range("I4").Value = 7.6
$ws.Range("J4").Value = 3.9
$ws.Range("K4").Value = 4.7
$ws.Range("P4").Value = 1.98
$ws.Range("Q4").Value = 1.8

# Row 5
$ws.Range("F5").Value = 1.79
$ws.Range("G5").Value = 1.97
$ws.Range("H5").Value = 4.5
$ws.Range("J5").Value = 3.3
$ws.Range("K5").Value = 3.9
$ws.Range("M5").Value = 1.09
$ws.Range("P5").Value = 1.68
$ws.Range("Q5").Value = 2.06
$ws.Range("T5").Value = 1.98
$ws.Range("U5").Value = 1.8
$ws.Range("AB5").Value = 970
$ws.Range("AF5").Value = 12
$ws.Range("AG5").Value = 12
$ws.Range("AI5").Value = 110

# Row 6
$ws.Range("P6").Value = 2.82
$ws.Range("Q6").Value = 1.3

# Row 7
$ws.Range("P7").Value = 2.74

# Row 8
$ws.Range("F8").Value = 2.68
$ws.Range("G8").Value = 3.1
$ws.Range("H8").Value = 3.3
$ws.Range("I8").Value = 3.95
$ws.Range("J8").Value = 2.7
$ws.Range("K8").Value = 2.96

# Row 9
$ws.Range("P9").Value = 1.59
$ws.Range("Q9").Value = 2.42

# Row 10
$ws.Range("F10").Value = 1.79
$ws.Range("I10").Value = 7.2
$ws.Range("J10").Value = 3.35
$ws.Range("K10").Value = 4.3
$ws.Range("P10").Value = 1.67
$ws.Range("Q10").Value = 2.22

# Row 12
$ws.Range("G12").Value = 2.3
$ws.Range("J12").Value = 3.3
$ws.Range("P12").Value = 1.8

# Row 14
$ws.Range("H14").Value = 5.6
$ws.Range("J14").Value = 3.75
